$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "48×12=" "61×30="
Replace-Text "92×35=" "98×60="
Replace-Text "79×87=" "32×63="
Replace-Text "74×71=" "85×96="
Replace-Text "95×89=" "88×23="
Replace-Text "79×41=" "15×33="
Replace-Text "60×88=" "84×74="
Replace-Text "55×95=" "20×19="
Replace-Text "25×71=" "18×17="
Replace-Text "43×44=" "76×19="
Replace-Text "28×38=" "43×66="
Replace-Text "74×36=" "23×25="
Replace-Text "99×84=" "64×49="
Replace-Text "87×24=" "27×93="
Replace-Text "58×82=" "59×41="
Replace-Text "32×33=" "85×87="
Replace-Text "37×20=" "40×47="
Replace-Text "89×14=" "49×35="
Replace-Text "20×78=" "31×27="
Replace-Text "60×52=" "42×78="
Replace-Text "16×14=" "65×37="
Replace-Text "15×38=" "75×79="
Replace-Text "52×66=" "81×73="
Replace-Text "87×77=" "75×39="
Replace-Text "51×55=" "39×24="
